$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("toybox_0_7_5_cppcheck")
$tbl = $ws2.ListObjects.Item(1)

# --- 1. Re-sort the table by filename (column B) ascending (was num_configs desc) ---
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws2.Range("B2:B21"))
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# --- 2. Fill in previously-blank classification (column H) values ---
$ws2.Range("H5").Value = $true    # bzcat.c
$ws2.Range("H7").Value = $false   # cp.c (line 343)
$ws2.Range("H8").Value = $false   # cp.c (line 409)
$ws2.Range("H11").Value = $false  # hwclock.c (line 96)
$ws2.Range("H12").Value = $false  # kill.c
$ws2.Range("H16").Value = $false  # main.c
$ws2.Range("H17").Value = $true   # netstat.c

# --- 3. Add two new table columns: manual_features, num_manual_features ---
$colI = $tbl.ListColumns.Add()
$ws2.Range("I1").Value = "manual_features"
$colJ = $tbl.ListColumns.Add()
$ws2.Range("J1").Value = "num_manual_features"

# --- 4. Populate the manual feature annotations for specific rows ---
$ws2.Range("I6").Value = "CONFIG_CHVT"
$ws2.Range("J6").Value = 1
$ws2.Range("I9").Value = "CONFIG_DATE"
$ws2.Range("J9").Value = 1
$ws2.Range("I10").Value = "CONFIG_HWCLOCK"
$ws2.Range("J10").Value = 1
$ws2.Range("I14").Value = "CONFIG_LOSETUP"
$ws2.Range("J14").Value = 1
$ws2.Range("I18").Value = "CONFIG_LOSETUP"
$ws2.Range("J18").Value = 1
$ws2.Range("I21").Value = "CONFIG_VMSTAT"
$ws2.Range("J21").Value = 1
